$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "55.049.43"
$ws.Range("E2").Value = "  -5.00%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "2.878.73"
$ws.Range("E3").Value = "  -5.58%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.48"
$ws.Range("E5").Value = "  -5.75%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.23"
$ws.Range("E6").Value = "  -8.10%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.414"
$ws.Range("E8").Value = "  -6.86%  "

$ws.Range("B9").Value = "Toncoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.08"
$ws.Range("E9").Value = "  -5.72%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  -8.76%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.343"
$ws.Range("E11").Value = "  -7.66%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "3.372.50"
$ws.Range("E12").Value = "  -5.68%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.124"
$ws.Range("E13").Value = "  -4.54%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.15"
$ws.Range("E14").Value = "  -6.19%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000155"
$ws.Range("E15").Value = "  -8.83%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "55.113.84"
$ws.Range("E16").Value = "  -4.90%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.91"
$ws.Range("E17").Value = "  -5.58%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.880.75"
$ws.Range("E18").Value = "  -5.46%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.25"
$ws.Range("E19").Value = "  -6.13%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.50"
$ws.Range("E20").Value = "  -7.30%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "308.57"
$ws.Range("E21").Value = "  -8.90%  "

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("E22").Value = "  +0.27%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.474"
$ws.Range("E24").Value = "  -5.64%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "61.58"
$ws.Range("E25").Value = "  -5.42%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.997.72"
$ws.Range("E26").Value = "  -5.49%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.158"
$ws.Range("E28").Value = "  -6.83%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0822"
$ws.Range("E29").Value = "  -13.82%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.22"
$ws.Range("E30").Value = "  -9.94%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.84"
$ws.Range("E31").Value = "  -9.87%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.72"
$ws.Range("E32").Value = "  -6.21%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.41"
$ws.Range("E33").Value = "  -7.62%  "

$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.10"
$ws.Range("E34").Value = "  -11.23%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "148.89"
$ws.Range("E35").Value = "  -5.13%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.29"
$ws.Range("E36").Value = "  -9.95%  "

$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.48"
$ws.Range("E37").Value = "  -7.59%  "

$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.08"
$ws.Range("E38").Value = "  -4.55%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.16"
$ws.Range("E39").Value = "  -10.70%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0641"
$ws.Range("E40").Value = "  -7.46%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.97"
$ws.Range("E42").Value = "  -4.86%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.59"
$ws.Range("E43").Value = "  -7.91%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.624"
$ws.Range("E44").Value = "  -6.04%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.063.69"
$ws.Range("E45").Value = "  -11.28%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.30"
$ws.Range("E46").Value = "  -10.21%  "

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.78"
$ws.Range("E47").Value = "  -4.89%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.893"
$ws.Range("E48").Value = "  -11.38%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0227"
$ws.Range("E49").Value = "  -6.27%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.31"
$ws.Range("E50").Value = "  -7.41%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0832"
$ws.Range("E51").Value = "  -7.56%  "
